$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - SignUp test case
$ws.Range("A2").Value = "TC_01.Verify SignUp with valid credentials"
$ws.Range("C2").Value = "testdata_Demo.xlsx,SignupPage"
$ws.Range("E2").Value = "'1"
$ws.Range("F2").Value = "'1"
$ws.Range("G2").Value = "Complete Sign Up flow and Add Address,Payment Method"

# Row 3 - Buy Token test case
$ws.Range("A3").Value = "TC_02.Verify Buy Token with New and Existing Debit and Credit Cards"
$ws.Range("C3").Value = "testdata_Demo.xlsx,buyTokens"
$ws.Range("E3").Value = "'1"
$ws.Range("F3").Value = "'1"
$ws.Range("G3").Value = "Add Credit Card and Do Buy Tokens"

# Row 4 - Withdraw Token test case
$ws.Range("A4").Value = "TC_03.Verify Withdraw Token with Gift Cards"
$ws.Range("C4").Value = "testdata_Demo.xlsx,buyTokens"
$ws.Range("E4").Value = "'1"
$ws.Range("F4").Value = "'1"
$ws.Range("G4").Value = "Do Withdraw Tokens with Gift Card"

# Row 5 - Add Card invalid data test case
$ws.Range("A5").Value = "TC_04.Verify Add Debit and Credit Card invalid data"
$ws.Range("C5").Value = "testdata_Demo.xlsx,PaymentMethods"
$ws.Range("E5").Value = "'1"
$ws.Range("F5").Value = "'1"
$ws.Range("G5").Value = "Add Card flow verify error messages"

# Row 6 - Scan Code / My QR Code test case
$ws.Range("A6").Value = "TC_05.Verify Scan Code and My QR Code"
$ws.Range("C6").Value = "testdata_Demo.xlsx,notifications"
$ws.Range("E6").Value = "'1"
$ws.Range("F6").Value = "'1"
$ws.Range("G6").Value = "Verify Scan Code"

# Row 7 - Send Transaction test case
$ws.Range("A7").Value = "TC_06.Verify Send Transaction"
$ws.Range("C7").Value = "testdata_Demo.xlsx,notifications"
$ws.Range("E7").Value = "'1"
$ws.Range("F7").Value = "'1"
$ws.Range("G7").Value = "Send Tokens from One User To another User"
$ws.Range("G7").Style = "Normal"

# Row 8 - Request Transaction test case
$ws.Range("A8").Value = "TC_07.Verify Request Transaction"
$ws.Range("C8").Value = "testdata_Demo.xlsx,notifications"
$ws.Range("E8").Value = "'1"
$ws.Range("F8").Value = "'1"
$ws.Range("G8").Value = "Request Tokens From one User to Another User"
$ws.Range("G8").Style = "Normal"

# Row 9 - Filters test case
$ws.Range("A9").Value = "TC_08.Verify Filters with one Transactions Type"
$ws.Range("C9").Value = "testdata_Demo.xlsx,filters"
$ws.Range("E9").Value = "'1"
$ws.Range("F9").Value = "'1"
$ws.Range("G9").Value = "Filters"
